$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rich-text header strings (Volume/Number and date range) ---
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "20"

$c9 = $ws.Range("C9")
$c9.Characters(27, 8).Text = "5/15/2023"
$c9.Characters(47, 9).Text = "5/21/2023"

# --- Update weekly crime statistics table (rows 14-30) ---
$ws.Range("C14").Value = 7
$ws.Range("D14").Value = 13
$ws.Range("E14").Value = -46.153846153846
$ws.Range("F14").Value = 24
$ws.Range("G14").Value = 43
$ws.Range("H14").Value = -44.186046511627
$ws.Range("I14").Value = 146
$ws.Range("J14").Value = 168
$ws.Range("K14").Value = -13.095238095238
$ws.Range("L14").Value = -19.337016574585
$ws.Range("M14").Value = -16.571428571428
$ws.Range("N14").Value = -79.862068965517
$ws.Range("D15").Value = 37
$ws.Range("E15").Value = -13.513513513513
$ws.Range("F15").Value = 115
$ws.Range("G15").Value = 115
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 573
$ws.Range("J15").Value = 608
$ws.Range("K15").Value = -5.756578947368
$ws.Range("L15").Value = 7.504690431519
$ws.Range("M15").Value = 19.624217118997
$ws.Range("N15").Value = -52.605459057072
$ws.Range("C16").Value = 305
$ws.Range("D16").Value = 327
$ws.Range("E16").Value = -6.727828746177
$ws.Range("F16").Value = 1213
$ws.Range("G16").Value = 1256
$ws.Range("H16").Value = -3.423566878980
$ws.Range("I16").Value = 5928
$ws.Range("J16").Value = 6075
$ws.Range("K16").Value = -2.419753086419
$ws.Range("L16").Value = 37.732342007434
$ws.Range("M16").Value = -13.181019332161
$ws.Range("N16").Value = -81.490617291660
$ws.Range("C17").Value = 514
$ws.Range("D17").Value = 554
$ws.Range("E17").Value = -7.220216606498
$ws.Range("F17").Value = 2043
$ws.Range("G17").Value = 2059
$ws.Range("H17").Value = -0.777076250607
$ws.Range("I17").Value = 9977
$ws.Range("J17").Value = 9261
$ws.Range("K17").Value = 7.731346506856
$ws.Range("L17").Value = 29.908854166666
$ws.Range("M17").Value = 61.361798479702
$ws.Range("N17").Value = -31.981183528770
$ws.Range("C18").Value = 271
$ws.Range("D18").Value = 274
$ws.Range("E18").Value = -1.094890510948
$ws.Range("F18").Value = 1026
$ws.Range("G18").Value = 1174
$ws.Range("H18").Value = -12.606473594548
$ws.Range("I18").Value = 5466
$ws.Range("J18").Value = 5917
$ws.Range("K18").Value = -7.622105796856
$ws.Range("L18").Value = 23.525423728813
$ws.Range("M18").Value = -18.588025022341
$ws.Range("N18").Value = -85.568697856162
$ws.Range("C19").Value = 997
$ws.Range("D19").Value = 980
$ws.Range("E19").Value = 1.734693877551
$ws.Range("F19").Value = 3780
$ws.Range("G19").Value = 3833
$ws.Range("H19").Value = -1.382728932950
$ws.Range("I19").Value = 18581
$ws.Range("J19").Value = 18804
$ws.Range("K19").Value = -1.185917889810
$ws.Range("L19").Value = 51.212565104166
$ws.Range("M19").Value = 37.882160878599
$ws.Range("N19").Value = -39.721005677210
$ws.Range("C20").Value = 292
$ws.Range("D20").Value = 247
$ws.Range("E20").Value = 18.218623481781
$ws.Range("F20").Value = 1237
$ws.Range("G20").Value = 953
$ws.Range("H20").Value = 29.800629590766
$ws.Range("I20").Value = 5723
$ws.Range("J20").Value = 4919
$ws.Range("K20").Value = 16.344785525513
$ws.Range("L20").Value = 80.195214105793
$ws.Range("M20").Value = 51.082365364308
$ws.Range("N20").Value = -86.868733221669
$ws.Range("C21").Value = 2418
$ws.Range("D21").Value = 2432
$ws.Range("E21").Value = -0.575657894736
$ws.Range("F21").Value = 9438
$ws.Range("G21").Value = 9433
$ws.Range("H21").Value = 0.053005406551
$ws.Range("I21").Value = 46394
$ws.Range("J21").Value = 45752
$ws.Range("K21").Value = 1.403217345689
$ws.Range("L21").Value = 42.369656611532
$ws.Range("M21").Value = 23.247350104933
$ws.Range("N21").Value = -71.168271053302
$ws.Range("C22").Value = 46
$ws.Range("D22").Value = 59
$ws.Range("E22").Value = -22.033898305084
$ws.Range("F22").Value = 172
$ws.Range("G22").Value = 202
$ws.Range("H22").Value = -14.851485148514
$ws.Range("I22").Value = 837
$ws.Range("J22").Value = 911
$ws.Range("K22").Value = -8.122941822173
$ws.Range("L22").Value = 48.141592920354
$ws.Range("M22").Value = 5.283018867924
$ws.Range("C23").Value = 101
$ws.Range("D23").Value = 110
$ws.Range("E23").Value = -8.181818181818
$ws.Range("F23").Value = 445
$ws.Range("G23").Value = 464
$ws.Range("H23").Value = -4.094827586206
$ws.Range("I23").Value = 2322
$ws.Range("J23").Value = 2165
$ws.Range("K23").Value = 7.251732101616
$ws.Range("L23").Value = 17.629179331307
$ws.Range("M23").Value = 59.917355371900
$ws.Range("C24").Value = 2126
$ws.Range("D24").Value = 2251
$ws.Range("E24").Value = -5.553087516659
$ws.Range("F24").Value = 8293
$ws.Range("G24").Value = 8853
$ws.Range("H24").Value = -6.325539365186
$ws.Range("I24").Value = 40887
$ws.Range("J24").Value = 41455
$ws.Range("K24").Value = -1.370160414907
$ws.Range("L24").Value = 40.009588056021
$ws.Range("M24").Value = 40.374909877433
$ws.Range("C25").Value = 884
$ws.Range("D25").Value = 846
$ws.Range("E25").Value = 4.491725768321
$ws.Range("F25").Value = 3601
$ws.Range("G25").Value = 3341
$ws.Range("H25").Value = 7.782101167315
$ws.Range("I25").Value = 16205
$ws.Range("J25").Value = 15287
$ws.Range("K25").Value = 6.005102374566
$ws.Range("L25").Value = 35.221962616822
$ws.Range("M25").Value = -4.524833559182
$ws.Range("C26").Value = 54
$ws.Range("D26").Value = 60
$ws.Range("E26").Value = -10
$ws.Range("F26").Value = 200
$ws.Range("G26").Value = 206
$ws.Range("H26").Value = -2.912621359223
$ws.Range("I26").Value = 935
$ws.Range("J26").Value = 1001
$ws.Range("K26").Value = -6.593406593406
$ws.Range("L26").Value = 4.469273743016
$ws.Range("C27").Value = 114
$ws.Range("D27").Value = 117
$ws.Range("E27").Value = -2.564102564102
$ws.Range("F27").Value = 452
$ws.Range("G27").Value = 428
$ws.Range("H27").Value = 5.607476635514
$ws.Range("I27").Value = 1975
$ws.Range("J27").Value = 1858
$ws.Range("K27").Value = 6.297093649085
$ws.Range("L27").Value = 19.263285024154
$ws.Range("C28").Value = 19
$ws.Range("D28").Value = 35
$ws.Range("E28").Value = -45.714285714285
$ws.Range("F28").Value = 80
$ws.Range("G28").Value = 118
$ws.Range("H28").Value = -32.203389830508
$ws.Range("I28").Value = 416
$ws.Range("J28").Value = 554
$ws.Range("K28").Value = -24.909747292418
$ws.Range("L28").Value = -28.151986183074
$ws.Range("M28").Value = -26.241134751773
$ws.Range("N28").Value = -80.533458118858
$ws.Range("C29").Value = 17
$ws.Range("D29").Value = 30
$ws.Range("E29").Value = -43.333333333333
$ws.Range("F29").Value = 71
$ws.Range("G29").Value = 99
$ws.Range("H29").Value = -28.282828282828
$ws.Range("I29").Value = 352
$ws.Range("J29").Value = 471
$ws.Range("K29").Value = -25.265392781316
$ws.Range("L29").Value = -31.115459882583
$ws.Range("M29").Value = -24.137931034482
$ws.Range("N29").Value = -81.789963786859
$ws.Range("C30").Value = 8
$ws.Range("D30").Value = 13
$ws.Range("E30").Value = -38.461538461538
$ws.Range("F30").Value = 37
$ws.Range("G30").Value = 56
$ws.Range("H30").Value = -33.928571428571
$ws.Range("I30").Value = 185
$ws.Range("J30").Value = 271
$ws.Range("K30").Value = -31.734317343173
$ws.Range("L30").Value = -7.035175879396
